# Added Sortable On FrontEnd
# ---------------------------------------------------------------------------
# Insert two new rows (37,38) right above the old "DotNet Migrations" header
# row (old row 38) on Sheet1, fill them with a new DotNet-migration command
# pair, then fix up the two hyperlinks that used to sit on B50/B51 (now
# B52/B53 after the shift) and restore the reported selection.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Push everything from row 38 down by two rows.
$ws.Rows("37:38").Insert()

# 2) Populate the two freshly inserted rows.
$ws.Cells.Item(37, 2).Value = "dotnet ef database update AddEventLanguages"
$ws.Cells.Item(37, 3).Value = "Swich To Migration"

$ws.Cells.Item(38, 2).Value = "dotnet ef migrations remove"
$ws.Cells.Item(38, 3).Value = "Remove last migration"

# 3) The two hyperlinks that lived on B50/B51 now live on B52/B53 - the
#    engine does not auto-shift Hyperlinks on row insert, so re-home them.
function Move-Hyperlink($ws, $oldAddr, $newAddr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $oldAddr) {
            $url = $h.Address
            $h.Delete()
            $ws.Hyperlinks.Add($ws.Range($newAddr), $url) | Out-Null
            return
        }
    }
}

Move-Hyperlink $ws '$B$50' "B52"
Move-Hyperlink $ws '$B$51' "B53"

# 4) Restore the reported view/selection state.
$ws.Activate()
try {
    $aw = $excel.ActiveWindow
    $aw.ScrollRow = 19
    $aw.ScrollColumn = 1
} catch {
}
$ws.Range("C38").Select()
